# Update timestamp embedded in test-data email addresses from
# 20251109_012452 -> 20251109_013943, on both sheets that contain them.

$wb = $excel.ActiveWorkbook

$oldStamp = "20251109_012452"
$newStamp = "20251109_013943"

# Sheet "UsuariosRegistro": column C (E-Mail), rows 2-6
$wsUsuarios = $wb.Worksheets.Item("UsuariosRegistro")
$wsUsuarios.Range("C2").Value = "juan.perez+$newStamp@test.com"
$wsUsuarios.Range("C3").Value = "maria.gonzalez+$newStamp@test.com"
$wsUsuarios.Range("C4").Value = "carlos.rodriguez+$newStamp@test.com"
$wsUsuarios.Range("C5").Value = "ana.martinez+$newStamp@test.com"
$wsUsuarios.Range("C6").Value = "luis.garcia+$newStamp@test.com"

# Sheet "LoginData": column A (Email), rows 2-3 reuse the first two addresses
$wsLogin = $wb.Worksheets.Item("LoginData")
$wsLogin.Range("A2").Value = "juan.perez+$newStamp@test.com"
$wsLogin.Range("A3").Value = "maria.gonzalez+$newStamp@test.com"
